$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1, reusing the
# same header formatting (style) as the existing header row by copying the
# format from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row (2-51).
for ($i = 2; $i -le 51; $i++) {
    $ws.Cells.Item($i, 30).Value = 88
    $ws.Cells.Item($i, 31).Value = 74
    $ws.Cells.Item($i, 32).Value = 0
}
